$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: bahn / ben
$ws.Range("A13").Value = "bahn"
$ws.Range("B13").Value = "ben"
$ws.Range("C13").Value = 0.5
$ws.Range("D13").Value = 0.545
$ws.Range("E13").Value = 0.522
$ws.Range("F13").Value = 0.846
$ws.Range("G13").Value = 0.75
$ws.Range("H13").Value = 0.795

# Row 14: auto-first / ben
$ws.Range("A14").Value = "auto-first"
$ws.Range("B14").Value = "ben"
$ws.Range("C14").Value = 0.14
$ws.Range("D14").Value = 0.212
$ws.Range("E14").Value = 0.169
$ws.Range("F14").Value = 0.34
$ws.Range("G14").Value = 0.386
$ws.Range("H14").Value = 0.362

# Row 15: auto-first / bahn
$ws.Range("A15").Value = "auto-first"
$ws.Range("B15").Value = "bahn"
$ws.Range("C15").Value = 0.14
$ws.Range("D15").Value = 0.194
$ws.Range("E15").Value = 0.163
$ws.Range("F15").Value = 0.32
$ws.Range("G15").Value = 0.41
$ws.Range("H15").Value = 0.36

# Row 16: averages of 14:15
$ws.Range("C16").Formula = "=AVERAGE(C14:C15)"
$ws.Range("D16").Formula = "=AVERAGE(D14:D15)"
$ws.Range("F16").Formula = "=AVERAGE(F14:F15)"
$ws.Range("G16").Formula = "=AVERAGE(G14:G15)"

# Row 18: auto-recent / ben
$ws.Range("A18").Value = "auto-recent"
$ws.Range("B18").Value = "ben"
$ws.Range("C18").Value = 0.327
$ws.Range("D18").Value = 0.485
$ws.Range("E18").Value = 0.39
$ws.Range("F18").Value = 0.551
$ws.Range("G18").Value = 0.614
$ws.Range("H18").Value = 0.581

# Row 19: auto-recent / bahn
$ws.Range("A19").Value = "auto-recent"
$ws.Range("B19").Value = "bahn"
$ws.Range("C19").Value = 0.286
$ws.Range("D19").Value = 0.389
$ws.Range("E19").Value = 0.329
$ws.Range("F19").Value = 0.469
$ws.Range("G19").Value = 0.59
$ws.Range("H19").Value = 0.523

# Row 20: averages of 18:19
$ws.Range("C20").Formula = "=AVERAGE(C18:C19)"
$ws.Range("D20").Formula = "=AVERAGE(D18:D19)"
$ws.Range("F20").Formula = "=AVERAGE(F18:F19)"
$ws.Range("G20").Formula = "=AVERAGE(G18:G19)"

# Row 22: auto-self / ben
$ws.Range("A22").Value = "auto-self"
$ws.Range("B22").Value = "ben"
$ws.Range("C22").Value = 0.293
$ws.Range("D22").Value = 0.364
$ws.Range("E22").Value = 0.324
$ws.Range("F22").Value = 0.488
$ws.Range("G22").Value = 0.455
$ws.Range("H22").Value = 0.471

# Row 23: auto-self / bahn
$ws.Range("A23").Value = "auto-self"
$ws.Range("B23").Value = "bahn"
$ws.Range("C23").Value = 0.317
$ws.Range("D23").Value = 0.361
$ws.Range("E23").Value = 0.338
$ws.Range("F23").Value = 0.488
$ws.Range("G23").Value = 0.513
$ws.Range("H23").Value = 0.5

# Row 24: averages of 22:23
$ws.Range("C24").Formula = "=AVERAGE(C22:C23)"
$ws.Range("D24").Formula = "=AVERAGE(D22:D23)"
$ws.Range("F24").Formula = "=AVERAGE(F22:F23)"
$ws.Range("G24").Formula = "=AVERAGE(G22:G23)"

# Row 26: auto-self2 / ben
$ws.Range("A26").Value = "auto-self2"
$ws.Range("B26").Value = "ben"
$ws.Range("C26").Value = 0.327
$ws.Range("D26").Value = 0.485
$ws.Range("E26").Value = 0.39
$ws.Range("F26").Value = 0.531
$ws.Range("G26").Value = 0.591
$ws.Range("H26").Value = 0.559

# Row 27: auto-self2 / bahn
$ws.Range("A27").Value = "auto-self2"
$ws.Range("B27").Value = "bahn"
$ws.Range("C27").Value = 0.327
$ws.Range("D27").Value = 0.444
$ws.Range("E27").Value = 0.376
$ws.Range("F27").Value = 0.51
$ws.Range("G27").Value = 0.641
$ws.Range("H27").Value = 0.568

# Row 28: averages of 26:27
$ws.Range("C28").Formula = "=AVERAGE(C26:C27)"
$ws.Range("D28").Formula = "=AVERAGE(D26:D27)"
$ws.Range("F28").Formula = "=AVERAGE(F26:F27)"
$ws.Range("G28").Formula = "=AVERAGE(G26:G27)"

# Update selection to match target state
$ws.Range("E28").Select()
